# v1.3 time was changed
# - "Rzeczywisty czas pracy" (actual time worked) for the "Sterowanie postacią"
#   milestone (row 5) changed from 30 to 90.
# - A note explaining the change was added next to it.
# - The view/selection was left on the newly edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Update the actual time worked for row 5 ("Sterowanie postacią")
$ws.Range("E5").Value = 90

# Add the explanatory note in the "Problemy, które wystąpiły/uwagi" column
$ws.Range("F5").Value = "Duże problemy z poprawieniem błędów w poruszaniu się. Więcej w post mortem"

# Leave the selection / view on the edited row
$ws.Range("F5").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
